$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "7.20", "0.999", "26.20") are preserved exactly instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.868.95"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "3.454.89"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "574.47"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "159.62"
$ws.Range("E6").Value = "  -1.84%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.454.50"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  -5.82%  "

$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("E11").Value = "  -2.83%  "

$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").Value = "4.052.35"
$ws.Range("E13").Value = "  -0.95%  "

$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("D15").Value = "27.65"
$ws.Range("E15").Value = "  -3.93%  "

$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -9.21%  "

$ws.Range("D17").Value = "64.965.12"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "3.440.50"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("E19").Value = "  -3.78%  "

$ws.Range("D20").Value = "13.78"
$ws.Range("E20").Value = "  -4.32%  "

$ws.Range("D21").Value = "379.39"
$ws.Range("E21").Value = "  -1.59%  "

$ws.Range("D22").Value = "7.95"
$ws.Range("E22").Value = "  -3.48%  "

$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "72.30"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("E25").Value = "  -3.20%  "

$ws.Range("D26").Value = "0.0000121"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "1.45"
$ws.Range("E30").Value = "  -4.44%  "

$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("E32").Value = "  -2.34%  "

$ws.Range("D33").Value = "23.20"
$ws.Range("E33").Value = "  -2.25%  "

$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("D36").Value = "161.13"
$ws.Range("E36").Value = "  -0.80%  "

$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("D38").Value = "2.899.94"
$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").Value = "0.0750"
$ws.Range("E39").Value = "  -3.81%  "

$ws.Range("D40").Value = "26.20"
$ws.Range("E40").Value = "  -3.05%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "4.53"
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "42.98"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "6.52"
$ws.Range("E43").Value = "  -4.76%  "

$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").Value = "26.19"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("E46").Value = "  -2.96%  "

$ws.Range("D47").Value = "2.36"
$ws.Range("E47").Value = "  +7.26%  "

$ws.Range("D48").Value = "321.67"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("E51").Value = "  -3.80%  "
